$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the cells we are about to rewrite stay formatted/stored as text,
# matching the inline-string cells already used throughout this sheet
# (column D = Price, column E = Volume(1h), column G = Hora).
$ws.Range("D2:D51").NumberFormat = "@"
$ws.Range("E2:E51").NumberFormat = "@"
$ws.Range("G2:G51").NumberFormat = "@"

$ws.Range("D2").Value = "260.95"
$ws.Range("E2").Value = "-0.26%"
$ws.Range("G2").Value = "3"
$ws.Range("D3").Value = "27.55"
$ws.Range("E3").Value = "-0.12%"
$ws.Range("G3").Value = "3"
$ws.Range("D4").Value = "4.710"
$ws.Range("E4").Value = "-0.29%"
$ws.Range("G4").Value = "3"
$ws.Range("D5").Value = "0.06222"
$ws.Range("E5").Value = "2.48%"
$ws.Range("G5").Value = "3"
$ws.Range("E6").Value = "0.36%"
$ws.Range("G6").Value = "3"
$ws.Range("D7").Value = "0.8504"
$ws.Range("E7").Value = "-1.43%"
$ws.Range("G7").Value = "3"
$ws.Range("D8").Value = "0.9106"
$ws.Range("E8").Value = "-1.25%"
$ws.Range("G8").Value = "3"
$ws.Range("E9").Value = "0.01%"
$ws.Range("G9").Value = "3"
$ws.Range("D10").Value = "0.04812"
$ws.Range("E10").Value = "-7.24%"
$ws.Range("G10").Value = "3"
$ws.Range("D11").Value = "0.07079"
$ws.Range("E11").Value = "-0.56%"
$ws.Range("G11").Value = "3"
$ws.Range("D12").Value = "0.03130"
$ws.Range("E12").Value = "1.87%"
$ws.Range("G12").Value = "3"
$ws.Range("D13").Value = "0.09044"
$ws.Range("E13").Value = "-0.56%"
$ws.Range("G13").Value = "3"
$ws.Range("D14").Value = "0.001537"
$ws.Range("E14").Value = "0.30%"
$ws.Range("G14").Value = "3"
$ws.Range("D15").Value = "0.0006162"
$ws.Range("E15").Value = "1.25%"
$ws.Range("G15").Value = "3"
$ws.Range("D16").Value = "0.005991"
$ws.Range("E16").Value = "-3.35%"
$ws.Range("G16").Value = "3"
$ws.Range("E17").Value = "-0.57%"
$ws.Range("G17").Value = "3"
$ws.Range("D18").Value = "3.167"
$ws.Range("E18").Value = "-0.10%"
$ws.Range("G18").Value = "3"
$ws.Range("G19").Value = "3"
$ws.Range("G20").Value = "3"
$ws.Range("E21").Value = "1.58%"
$ws.Range("G21").Value = "3"
$ws.Range("E22").Value = "-0.22%"
$ws.Range("G22").Value = "3"
$ws.Range("D23").Value = "0.04251"
$ws.Range("E23").Value = "-0.07%"
$ws.Range("G23").Value = "3"
$ws.Range("D24").Value = "0.001219"
$ws.Range("E24").Value = "0.18%"
$ws.Range("G24").Value = "3"
$ws.Range("D25").Value = "0.004080"
$ws.Range("E25").Value = "4.25%"
$ws.Range("G25").Value = "3"
$ws.Range("E26").Value = "0.06%"
$ws.Range("G26").Value = "3"
$ws.Range("D27").Value = "0.0001640"
$ws.Range("E27").Value = "4.90%"
$ws.Range("G27").Value = "3"
$ws.Range("G28").Value = "3"
$ws.Range("G29").Value = "3"
$ws.Range("G30").Value = "3"
$ws.Range("G31").Value = "3"
$ws.Range("G32").Value = "3"
$ws.Range("G33").Value = "3"
$ws.Range("G34").Value = "3"
$ws.Range("G35").Value = "3"
$ws.Range("G36").Value = "3"
$ws.Range("G37").Value = "3"
$ws.Range("G38").Value = "3"
$ws.Range("G39").Value = "3"
$ws.Range("D40").Value = "0.03874"
$ws.Range("E40").Value = "-0.15%"
$ws.Range("G40").Value = "3"
$ws.Range("E41").Value = "-0.37%"
$ws.Range("G41").Value = "3"
$ws.Range("D42").Value = "0.004112"
$ws.Range("E42").Value = "-0.75%"
$ws.Range("G42").Value = "3"
$ws.Range("E43").Value = "-2.75%"
$ws.Range("G43").Value = "3"
$ws.Range("D44").Value = "0.01313"
$ws.Range("E44").Value = "-12.54%"
$ws.Range("G44").Value = "3"
$ws.Range("D45").Value = "0.00005136"
$ws.Range("E45").Value = "-2.85%"
$ws.Range("G45").Value = "3"
$ws.Range("E46").Value = "0.06%"
$ws.Range("G46").Value = "3"
$ws.Range("D47").Value = "0.03403"
$ws.Range("E47").Value = "-37.64%"
$ws.Range("G47").Value = "3"
$ws.Range("D48").Value = "0.06960"
$ws.Range("E48").Value = "-48.56%"
$ws.Range("G48").Value = "3"
$ws.Range("E49").Value = "0.06%"
$ws.Range("G49").Value = "3"
$ws.Range("E50").Value = "0.06%"
$ws.Range("G50").Value = "3"
$ws.Range("G51").Value = "3"
